$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder / relabel the header row (row 1) so columns are grouped by category
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "kitchens_2"
$ws.Range("C1").Value = "bedrooms_1"
$ws.Range("D1").Value = "bedrooms_2"
$ws.Range("E1").Value = "living_rooms_1"
$ws.Range("F1").Value = "living_rooms_2"

# Update the data rows (2-6) so that the "1" marker follows the same
# category label into its new column position. Row 7 is unchanged.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0

$ws.Range("A4").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1

$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
